$wb = $excel.ActiveWorkbook

$ws_LP1912 = $wb.Worksheets.Item("LP1912")
# Rows 2-3
$arr = New-Object 'object[,]' 2,5
$arr[0,0] = 'Última actualización: 20:32:02'
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[1,0] = 'Total filas: 144'
$arr[1,1] = $null
$arr[1,2] = $null
$arr[1,3] = $null
$arr[1,4] = $null
$ws_LP1912.Range("A2:E3").Value = $arr

# Rows 39-40
$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '16:52:37'
$arr[0,1] = '17:51'
$arr[0,2] = '215_EL PELIGRO'
$arr[0,3] = 59
$arr[0,4] = 'LP1912'
$arr[1,0] = '17:47:22'
$arr[1,1] = '17:51'
$arr[1,2] = '215B_EL PATO'
$arr[1,3] = 4
$arr[1,4] = 'LP1912'
$ws_LP1912.Range("A39:E40").Value = $arr

# Rows 48-49
$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '18:10:41'
$arr[0,1] = '18:11'
$arr[0,2] = '16_SANTA ANA'
$arr[0,3] = 1
$arr[0,4] = 'LP1912'
$arr[1,0] = '18:10:41'
$arr[1,1] = '18:11'
$arr[1,2] = '10_OLMOS'
$arr[1,3] = 1
$arr[1,4] = 'LP1912'
$ws_LP1912.Range("A48:E49").Value = $arr

# Rows 76-77
$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '17:47:22'
$arr[0,1] = '19:17'
$arr[0,2] = '27_EL RETIRO'
$arr[0,3] = 90
$arr[0,4] = 'LP1912'
$arr[1,0] = '18:44:34'
$arr[1,1] = '19:17'
$arr[1,2] = '16_SANTA ANA'
$arr[1,3] = 33
$arr[1,4] = 'LP1912'
$ws_LP1912.Range("A76:E77").Value = $arr

# Rows 93-94
$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '17:47:22'
$arr[0,1] = '19:40'
$arr[0,2] = '17X38_ROMERO'
$arr[0,3] = 113
$arr[0,4] = 'LP1912'
$arr[1,0] = '19:11:59'
$arr[1,1] = '19:40'
$arr[1,2] = '16_SANTA ANA'
$arr[1,3] = 29
$arr[1,4] = 'LP1912'
$ws_LP1912.Range("A93:E94").Value = $arr

# Rows 112-113
$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '20:11:56'
$arr[0,1] = '20:12'
$arr[0,2] = '11_ETCHEVERRY'
$arr[0,3] = 1
$arr[0,4] = 'LP1912'
$arr[1,0] = '20:11:56'
$arr[1,1] = '20:12'
$arr[1,2] = '16_SANTA ANA'
$arr[1,3] = 1
$arr[1,4] = 'LP1912'
$ws_LP1912.Range("A112:E113").Value = $arr

# Rows 120-149
$arr = New-Object 'object[,]' 30,5
$arr[0,0] = '20:32:02'
$arr[0,1] = '20:32'
$arr[0,2] = '10_OLMOS'
$arr[0,3] = 0
$arr[0,4] = 'LP1912'
$arr[1,0] = '20:32:02'
$arr[1,1] = '20:33'
$arr[1,2] = '16_SANTA ANA'
$arr[1,3] = 1
$arr[1,4] = 'LP1912'
$arr[2,0] = '20:11:56'
$arr[2,1] = '20:35'
$arr[2,2] = '16_SANTA ANA'
$arr[2,3] = 24
$arr[2,4] = 'LP1912'
$arr[3,0] = '20:32:02'
$arr[3,1] = '20:36'
$arr[3,2] = '16_SANTA ANA'
$arr[3,3] = 4
$arr[3,4] = 'LP1912'
$arr[4,0] = '19:11:59'
$arr[4,1] = '20:43'
$arr[4,2] = '215B_EL PATO'
$arr[4,3] = 92
$arr[4,4] = 'LP1912'
$arr[5,0] = '19:11:59'
$arr[5,1] = '20:44'
$arr[5,2] = '17X38_ROMERO'
$arr[5,3] = 93
$arr[5,4] = 'LP1912'
$arr[6,0] = '18:52:04'
$arr[6,1] = '20:44'
$arr[6,2] = '215B_EL PATO'
$arr[6,3] = 112
$arr[6,4] = 'LP1912'
$arr[7,0] = '18:52:04'
$arr[7,1] = '20:45'
$arr[7,2] = '17X38_ROMERO'
$arr[7,3] = 113
$arr[7,4] = 'LP1912'
$arr[8,0] = '20:11:56'
$arr[8,1] = '20:49'
$arr[8,2] = '23_HERNANDEZ'
$arr[8,3] = 38
$arr[8,4] = 'LP1912'
$arr[9,0] = '19:54:49'
$arr[9,1] = '20:50'
$arr[9,2] = '23_HERNANDEZ'
$arr[9,3] = 56
$arr[9,4] = 'LP1912'
$arr[10,0] = '20:32:02'
$arr[10,1] = '20:51'
$arr[10,2] = '23_HERNANDEZ'
$arr[10,3] = 19
$arr[10,4] = 'LP1912'
$arr[11,0] = '19:35:31'
$arr[11,1] = '20:52'
$arr[11,2] = '23_HERNANDEZ'
$arr[11,3] = 77
$arr[11,4] = 'LP1912'
$arr[12,0] = '19:54:49'
$arr[12,1] = '20:56'
$arr[12,2] = '27_EL RETIRO'
$arr[12,3] = 62
$arr[12,4] = 'LP1912'
$arr[13,0] = '20:32:02'
$arr[13,1] = '20:57'
$arr[13,2] = '27_EL RETIRO'
$arr[13,3] = 25
$arr[13,4] = 'LP1912'
$arr[14,0] = '19:11:59'
$arr[14,1] = '21:01'
$arr[14,2] = '215A_EL PATO'
$arr[14,3] = 110
$arr[14,4] = 'LP1912'
$arr[15,0] = '20:32:02'
$arr[15,1] = '21:02'
$arr[15,2] = '215A_EL PATO'
$arr[15,3] = 30
$arr[15,4] = 'LP1912'
$arr[16,0] = '19:11:59'
$arr[16,1] = '21:02'
$arr[16,2] = '27_EL RETIRO'
$arr[16,3] = 111
$arr[16,4] = 'LP1912'
$arr[17,0] = '19:47:58'
$arr[17,1] = '21:06'
$arr[17,2] = '27_EL RETIRO'
$arr[17,3] = 79
$arr[17,4] = 'LP1912'
$arr[18,0] = '19:35:31'
$arr[18,1] = '21:10'
$arr[18,2] = '27_EL RETIRO'
$arr[18,3] = 95
$arr[18,4] = 'LP1912'
$arr[19,0] = '19:35:31'
$arr[19,1] = '21:23'
$arr[19,2] = '10_OLMOS'
$arr[19,3] = 108
$arr[19,4] = 'LP1912'
$arr[20,0] = '20:32:02'
$arr[20,1] = '21:24'
$arr[20,2] = '10_OLMOS'
$arr[20,3] = 52
$arr[20,4] = 'LP1912'
$arr[21,0] = '20:11:56'
$arr[21,1] = '21:34'
$arr[21,2] = '23_HERNANDEZ'
$arr[21,3] = 83
$arr[21,4] = 'LP1912'
$arr[22,0] = '20:11:56'
$arr[22,1] = '21:48'
$arr[22,2] = '11_ETCHEVERRY'
$arr[22,3] = 97
$arr[22,4] = 'LP1912'
$arr[23,0] = '19:54:49'
$arr[23,1] = '21:49'
$arr[23,2] = '11_ETCHEVERRY'
$arr[23,3] = 115
$arr[23,4] = 'LP1912'
$arr[24,0] = '20:11:56'
$arr[24,1] = '21:55'
$arr[24,2] = '84_COLONIA URQUIZA-ESC 49'
$arr[24,3] = 104
$arr[24,4] = 'LP1912'
$arr[25,0] = '20:32:02'
$arr[25,1] = '21:56'
$arr[25,2] = '84_COLONIA URQUIZA-ESC 49'
$arr[25,3] = 84
$arr[25,4] = 'LP1912'
$arr[26,0] = '20:32:02'
$arr[26,1] = '22:19'
$arr[26,2] = '10_OLMOS'
$arr[26,3] = 107
$arr[26,4] = 'LP1912'
$arr[27,0] = '20:32:02'
$arr[27,1] = '22:21'
$arr[27,2] = '23_HERNANDEZ'
$arr[27,3] = 109
$arr[27,4] = 'LP1912'
$arr[28,0] = '20:32:02'
$arr[28,1] = '22:26'
$arr[28,2] = '15_ABASTO'
$arr[28,3] = 114
$arr[28,4] = 'LP1912'
$arr[29,0] = '20:32:02'
$arr[29,1] = '22:31'
$arr[29,2] = '215C_EL PATO'
$arr[29,3] = 119
$arr[29,4] = 'LP1912'
$ws_LP1912.Range("A120:E149").Value = $arr


$ws_LP1912_215 = $wb.Worksheets.Item("LP1912-215")
# Rows 2-3
$arr = New-Object 'object[,]' 2,5
$arr[0,0] = 'Última actualización: 20:32:02'
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[1,0] = 'Total filas: 17'
$arr[1,1] = $null
$arr[1,2] = $null
$arr[1,3] = $null
$arr[1,4] = $null
$ws_LP1912_215.Range("A2:E3").Value = $arr

# Rows 21-22
$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '20:32:02'
$arr[0,1] = '21:02'
$arr[0,2] = '215A_EL PATO'
$arr[0,3] = 30
$arr[0,4] = 'LP1912'
$arr[1,0] = '20:32:02'
$arr[1,1] = '22:31'
$arr[1,2] = '215C_EL PATO'
$arr[1,3] = 119
$arr[1,4] = 'LP1912'
$ws_LP1912_215.Range("A21:E22").Value = $arr


$ws_6203_6173 = $wb.Worksheets.Item("6203-6173")
# Rows 2-3
$arr = New-Object 'object[,]' 2,5
$arr[0,0] = 'Última actualización: 20:32:02'
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[1,0] = 'Total filas: 17'
$arr[1,1] = $null
$arr[1,2] = $null
$arr[1,3] = $null
$arr[1,4] = $null
$ws_6203_6173.Range("A2:E3").Value = $arr

# Rows 21-22
$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '20:32:02'
$arr[0,1] = '22:13'
$arr[0,2] = '215B_LP-P MOR-1 Y 57'
$arr[0,3] = 101
$arr[0,4] = 'L6173'
$arr[1,0] = '20:32:02'
$arr[1,1] = '22:19'
$arr[1,2] = '215A_LA PLATA'
$arr[1,3] = 107
$arr[1,4] = 'L6173'
$ws_6203_6173.Range("A21:E22").Value = $arr

